# New weekly price record is inserted right after the current row 75,
# shifting the existing historical rows (76-85) down to (77-86).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 76; rows 76-85 shift down to 77-86.
$ws.Rows.Item(76).Insert()

# Populate the newly inserted row 76 with the new weekly record.
$ws.Cells.Item(76,1).Value = 8
$ws.Cells.Item(76,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(76,3).Value = "Coquimbo"
$ws.Cells.Item(76,4).Value = 44474
$ws.Cells.Item(76,5).Value = 4
$ws.Cells.Item(76,6).Value = 100112044
$ws.Cells.Item(76,7).Value = "Perejil"
$ws.Cells.Item(76,8).Value = "Sin especificar"
$ws.Cells.Item(76,9).Value = "Primera"
$ws.Cells.Item(76,10).Value = 2800
$ws.Cells.Item(76,11).Value = 1500
$ws.Cells.Item(76,12).Value = 2000
$ws.Cells.Item(76,13).Value = 1750
$ws.Cells.Item(76,14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(76,15).Value = "Provincia del Elquí"
$ws.Cells.Item(76,16).Value = 1167
$ws.Cells.Item(76,17).Value = 1.5
$ws.Cells.Item(76,18).Value = "Hortaliza"

Write-Output "Inserted new weekly record at row 76"
